# Refresh the "cryptos" price/volume snapshot (GitHub Actions style update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their text representation even though many
# values look numeric (e.g. "36.21"). Setting the NumberFormat to Text ("@") for
# the whole data range before assigning values prevents Excel from silently
# re-interpreting them as numbers/dates.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Price (column D) updates ---
$ws.Range("D2").Value = "43.187.15"
$ws.Range("D3").Value = "2.322.26"
$ws.Range("D5").Value = "302.75"
$ws.Range("D6").Value = "99.51"
$ws.Range("D7").Value = "0.508"
$ws.Range("D9").Value = "0.517"
$ws.Range("D10").Value = "36.21"
$ws.Range("D13").Value = "17.61"
$ws.Range("D15").Value = "2.682.77"
$ws.Range("D16").Value = "2.370.84"
$ws.Range("D17").Value = "0.798"
$ws.Range("D18").Value = "43.110.40"
$ws.Range("D19").Value = "13.11"
$ws.Range("D20").Value = "6.25"
$ws.Range("D21").Value = "0.0₃0910"
$ws.Range("D22").Value = "68.12"
$ws.Range("D23").Value = "241.70"
$ws.Range("D25").Value = "2.46"
$ws.Range("D27").Value = "25.47"
$ws.Range("D28").Value = "168.60"
$ws.Range("D29").Value = "34.25"
$ws.Range("D30").Value = "9.21"
$ws.Range("D32").Value = "5.20"
$ws.Range("D33").Value = "0.999"
$ws.Range("D34").Value = "4.74"
$ws.Range("D35").Value = "17.87"
$ws.Range("D37").Value = "0.0699"
$ws.Range("D42").Value = "1.999.96"
$ws.Range("D45").Value = "10.12"
$ws.Range("D46").Value = "17.62"
$ws.Range("D48").Value = "77.04"
$ws.Range("D49").Value = "55.00"
$ws.Range("D50").Value = "2.548.03"

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +1.84%  "
$ws.Range("E10").Value = "  +5.14%  "
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("E13").Value = "  -1.42%  "
$ws.Range("E14").Value = "  +1.86%  "
$ws.Range("E15").Value = "  +0.72%  "
$ws.Range("E16").Value = "  +5.13%  "
$ws.Range("E17").Value = "  -1.26%  "
$ws.Range("E19").Value = "  +6.34%  "
$ws.Range("E20").Value = "  +2.42%  "
$ws.Range("E21").Value = "  +0.52%  "
$ws.Range("E22").Value = "  +0.45%  "
$ws.Range("E23").Value = "  +2.11%  "
$ws.Range("E24").Value = "  -1.03%  "
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("E27").Value = "  +3.89%  "
$ws.Range("E28").Value = "  +0.35%  "
$ws.Range("E29").Value = "  +1.39%  "
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("E31").Value = "  -2.27%  "
$ws.Range("E32").Value = "  +3.67%  "
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("E34").Value = "  +3.70%  "
$ws.Range("E35").Value = "  +5.31%  "
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("E39").Value = "  +0.75%  "
$ws.Range("E40").Value = "  -1.54%  "
$ws.Range("E41").Value = "  +0.37%  "
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("E43").Value = "  +1.30%  "
$ws.Range("E44").Value = "  -4.92%  "
$ws.Range("E45").Value = "  +0.41%  "
$ws.Range("E46").Value = "  -0.66%  "
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("E48").Value = "  +9.68%  "
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("E50").Value = "  +0.80%  "
$ws.Range("E51").Value = "  +1.73%  "
